$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that may look like plain numbers need to be forced to text so Excel
# does not auto-convert them to numeric values (the source data are text labels).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "72.921.49"
Set-TextValue $ws.Cells.Item(2, 5) "  +3.49%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "3.978.02"
Set-TextValue $ws.Cells.Item(3, 5) "  +1.18%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  +0.07%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "589.14"
Set-TextValue $ws.Cells.Item(5, 5) "  +10.15%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "159.28"
Set-TextValue $ws.Cells.Item(6, 5) "  +8.57%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.686"
Set-TextValue $ws.Cells.Item(7, 5) "  +0.36%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.999"
Set-TextValue $ws.Cells.Item(8, 5) "  -0.18%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.749"
Set-TextValue $ws.Cells.Item(9, 5) "  +2.20%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 5) "  +2.24%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "54.28"
Set-TextValue $ws.Cells.Item(11, 5) "  -0.52%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 5) "  +0.83%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "10.88"
Set-TextValue $ws.Cells.Item(13, 5) "  +3.49%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "4.607.92"
Set-TextValue $ws.Cells.Item(14, 5) "  +1.29%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "3.984.40"
Set-TextValue $ws.Cells.Item(15, 5) "  +1.69%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 5) "  +8.95%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "14.03"
Set-TextValue $ws.Cells.Item(17, 5) "  +2.79%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "20.34"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.44%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 5) "  +0.46%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "72.591.29"
Set-TextValue $ws.Cells.Item(20, 5) "  +3.40%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "434.42"
Set-TextValue $ws.Cells.Item(21, 5) "  +3.03%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 5) "  +13.49%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "96.05"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.09%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "3.42"
Set-TextValue $ws.Cells.Item(24, 5) "  -2.61%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "14.31"
Set-TextValue $ws.Cells.Item(25, 5) "  +1.56%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "4.42"
Set-TextValue $ws.Cells.Item(26, 5) "  +21.59%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) "  -1.30%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 2) "Filecoin"
Set-TextValue $ws.Cells.Item(28, 3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(28, 4) "10.58"
Set-TextValue $ws.Cells.Item(28, 5) "  +0.85%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 2) "LEO"
Set-TextValue $ws.Cells.Item(29, 3) "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Cells.Item(29, 4) "5.95"
Set-TextValue $ws.Cells.Item(29, 5) "  +1.80%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "36.37"
Set-TextValue $ws.Cells.Item(30, 5) "  +0.57%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "7.84"
Set-TextValue $ws.Cells.Item(31, 5) "  +2.88%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 5) "  +3.76%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 5) "  +2.61%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "681.80"
Set-TextValue $ws.Cells.Item(34, 5) "  +0.74%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "48.42"
Set-TextValue $ws.Cells.Item(35, 5) "  -2.18%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 5) "  +9.54%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "0.0₃0874"
Set-TextValue $ws.Cells.Item(37, 5) "  +7.13%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 5) "  +0.53%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "3.40"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.74%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 5) "  -1.26%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "0.998"
Set-TextValue $ws.Cells.Item(41, 5) "  -0.13%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "3.33"
Set-TextValue $ws.Cells.Item(42, 5) "  +4.07%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "1.00"
Set-TextValue $ws.Cells.Item(43, 5) "  +0.45%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "10.91"
Set-TextValue $ws.Cells.Item(44, 5) "  +13.71%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 5) "  +1.79%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.150"
Set-TextValue $ws.Cells.Item(46, 5) "  +1.54%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 5) "  -1.70%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "3.42"
Set-TextValue $ws.Cells.Item(48, 5) "  +2.21%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "3.01"
Set-TextValue $ws.Cells.Item(49, 5) "  +1.84%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 5) "  +5.39%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "2.807.02"
Set-TextValue $ws.Cells.Item(51, 5) "  +12.06%  "
